$d = $word.ActiveDocument

# Locate the existing bullet paragraph that ends with
# "The server application that stores ... command line application"
# (it currently carries the _GoBack bookmark at its end).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "The server application that stores the training results*") {
        $target = $candidate
        break
    }
}

# ---------------------------------------------------------------
# 1) New paragraph: "Diagrams" (top level bullet, ListLevelNumber=1)
# ---------------------------------------------------------------
$target.Range.InsertParagraphAfter()
$pDiagrams = $target.Next()
$pDiagrams.Range.ListFormat.ListLevelNumber = 1
$pDiagrams.Range.Text = "Diagrams"

# ---------------------------------------------------------------
# 2) New paragraph: "Web app: use case diagram, structure diagram, class diagram"
#    (second level bullet, ListLevelNumber=2)
# ---------------------------------------------------------------
$pDiagrams.Range.InsertParagraphAfter()
$pWebApp = $pDiagrams.Next()
$pWebApp.Range.ListFormat.ListLevelNumber = 2
$pWebApp.Range.Text = "Web app: use case diagram, structure diagram, class diagram"

# ---------------------------------------------------------------
# 3) New paragraph: "Console app: use case diagram, flow diagram, layers diagram, class diagram"
#    split across 6 separate runs (matching the original authoring)
# ---------------------------------------------------------------
$pWebApp.Range.InsertParagraphAfter()
$pConsole = $pWebApp.Next()
$pConsole.Range.ListFormat.ListLevelNumber = 2
$consoleSegments = @(
    "Console app: ",
    "use case diagram, ",
    "flow",
    " diagram,",
    " layers diagram,",
    " class diagram"
)
$consoleFullText = [string]::Join("", $consoleSegments)
$pConsole.Range.Text = $consoleFullText

$consoleStart = $pConsole.Range.Start
$consoleEnd = $pConsole.Range.End - 1
$pos = $consoleStart
for ($i = 0; $i -lt $consoleSegments.Length - 1; $i++) {
    $pos = $pos + $consoleSegments[$i].Length
    $tail = $d.Range($pos, $consoleEnd)
    $tail.Font.Bold = 1
    $tail2 = $d.Range($pos, $consoleEnd)
    $tail2.Font.Bold = 0
}

# ---------------------------------------------------------------
# 4) New paragraph: "Server app: use case diagram, package diagram, class diagram"
#    with the _GoBack bookmark relocated between "package diagram," and " class diagram"
# ---------------------------------------------------------------
$pConsole.Range.InsertParagraphAfter()
$pServer = $pConsole.Next()
$pServer.Range.ListFormat.ListLevelNumber = 2
$serverLead = "Server app: use case diagram, package diagram,"
$serverTail = " class diagram"
$pServer.Range.Text = $serverLead + $serverTail

$serverStart = $pServer.Range.Start
$bookmarkPos = $serverStart + $serverLead.Length
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
